$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D, shifting existing D:K to E:L
$ws.Columns("D").Insert()

# Copy number formatting from column E (the old D, now shifted) into new column D
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the new period/values
$ws.Range("D7").Value2 = 43465
$ws.Range("D8").Value2 = 1072300
$ws.Range("D9").Value2 = 370300
$ws.Range("D10").Value2 = 702000
$ws.Range("D12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("D14").Value2 = 1800
$ws.Range("D15").Value2 = 197700
$ws.Range("D17").Value2 = 804400
$ws.Range("D18").Value2 = 267900
$ws.Range("D20").Value2 = 4500
$ws.Range("D21").Value2 = 470100
$ws.Range("D22").Value2 = 60400
$ws.Range("D23").Value2 = 212000
$ws.Range("D24").Value2 = 47200
$ws.Range("D25").Value2 = 0
$ws.Range("D26").Value2 = 164800
$ws.Range("D27").Value2 = 164800
$ws.Range("D28").Value2 = 0
$ws.Range("D29").Value2 = 0
$ws.Range("D30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("D32").Value2 = -4500
$ws.Range("D33").Value2 = 164800
$ws.Range("D34").Value2 = 0
$ws.Range("D35").Value2 = 164800
$ws.Range("D38").Value2 = 43465
$ws.Range("D41").Value2 = 264100
$ws.Range("D42").Value2 = 0
$ws.Range("D43").Value2 = 40700
$ws.Range("D44").Value2 = 0
$ws.Range("D45").Value2 = 13100
$ws.Range("D46").Value2 = 317900
$ws.Range("D47").Value2 = 0
$ws.Range("D48").Value2 = 848000
$ws.Range("D49").Value2 = 1126000
$ws.Range("D50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("D52").Value2 = 11400
$ws.Range("D53").Value2 = 0
$ws.Range("D54").Value2 = 2303200
$ws.Range("D57").Value2 = 20800
$ws.Range("D58").Value2 = 25300
$ws.Range("D59").Value2 = 87600
$ws.Range("D60").Value2 = 133700
$ws.Range("D61").Value2 = 1142100
$ws.Range("D62").Value2 = 252100
$ws.Range("D63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("D66").Value2 = 1527900
$ws.Range("D68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("D72").Value2 = 850300
$ws.Range("D73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("D76").Value2 = 775400
$ws.Range("D77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("D81").Value2 = 164800
$ws.Range("D83").Value2 = 197700
$ws.Range("D84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("D89").Value2 = 407800
$ws.Range("D91").Value2 = -217800
$ws.Range("D92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("D94").Value2 = -214300
$ws.Range("D96").Value2 = -42900
$ws.Range("D97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("D100").Value2 = -91100
$ws.Range("D101").Value2 = 0
$ws.Range("D102").Value2 = 102400

# Fix up shifted values in rows where the restated prior-period figures differ from a pure shift
$ws.Range("F8").Value2 = 819300
$ws.Range("E10").Value2 = 622900
$ws.Range("F10").Value2 = 522800
$ws.Range("E14").Value2 = 5900
$ws.Range("E17").Value2 = 723600
$ws.Range("F17").Value2 = 632300
$ws.Range("E18").Value2 = 236300
$ws.Range("F18").Value2 = 187100
$ws.Range("E21").Value2 = 418600
$ws.Range("F21").Value2 = 340100
$ws.Range("E23").Value2 = 190100
$ws.Range("F23").Value2 = 162000
$ws.Range("E24").Value2 = 68900
$ws.Range("F24").Value2 = 61700
$ws.Range("E26").Value2 = 121200
$ws.Range("F26").Value2 = 100300
$ws.Range("E27").Value2 = 121200
$ws.Range("F27").Value2 = 100300
$ws.Range("E29").Value2 = 114000
$ws.Range("E33").Value2 = 235200
$ws.Range("F33").Value2 = 100300
$ws.Range("E35").Value2 = 235200
$ws.Range("F35").Value2 = 100300
$ws.Range("E43").Value2 = 128600
$ws.Range("E59").Value2 = 102200
$ws.Range("E81").Value2 = 235200
$ws.Range("F81").Value2 = 100300
